$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (the source data is inline string text).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D39", "D42", "D43", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values row by row.
$ws.Range("D2").Value = "42.393.52"
$ws.Range("E2").Value = "  -2.36%  "

$ws.Range("D3").Value = "2.220.65"
$ws.Range("E3").Value = "  -2.09%  "

$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").Value = "108.15"
$ws.Range("E5").Value = "  -9.28%  "

$ws.Range("D6").Value = "295.75"
$ws.Range("E6").Value = "  +11.54%  "

$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -3.29%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  -3.13%  "

$ws.Range("D10").Value = "43.52"
$ws.Range("E10").Value = "  -8.24%  "

$ws.Range("D11").Value = "0.0913"
$ws.Range("E11").Value = "  -3.13%  "

$ws.Range("D12").Value = "54.38"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").Value = "8.76"
$ws.Range("E13").Value = "  -4.48%  "

$ws.Range("D14").Value = "0.996"
$ws.Range("E14").Value = "  +10.40%  "

$ws.Range("E15").Value = "  -2.50%  "

$ws.Range("D16").Value = "15.07"
$ws.Range("E16").Value = "  -2.16%  "

$ws.Range("D17").Value = "2.550.08"
$ws.Range("E17").Value = "  -2.26%  "

$ws.Range("D18").Value = "2.217.55"
$ws.Range("E18").Value = "  -2.28%  "

$ws.Range("D19").Value = "42.290.33"
$ws.Range("E19").Value = "  -2.85%  "

$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  +7.38%  "

$ws.Range("E21").Value = "  -4.06%  "

$ws.Range("D22").Value = "72.14"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "3.45"
$ws.Range("E23").Value = "  +20.41%  "

$ws.Range("D24").Value = "2.32"
$ws.Range("E24").Value = "  -2.75%  "

$ws.Range("D25").Value = "228.41"
$ws.Range("E25").Value = "  -2.93%  "

$ws.Range("D26").Value = "9.07"
$ws.Range("E26").Value = "  -4.53%  "

$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -1.75%  "

$ws.Range("D28").Value = "11.63"
$ws.Range("E28").Value = "  -3.04%  "

$ws.Range("E29").Value = "  -1.14%  "

$ws.Range("D30").Value = "38.17"
$ws.Range("E30").Value = "  -8.73%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "173.94"
$ws.Range("E31").Value = "  +1.25%  "

$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  -5.47%  "

$ws.Range("D33").Value = "20.97"
$ws.Range("E33").Value = "  -2.94%  "

$ws.Range("D34").Value = "0.0899"
$ws.Range("E34").Value = "  -1.63%  "

$ws.Range("D35").Value = "5.58"
$ws.Range("E35").Value = "  -2.22%  "

$ws.Range("D36").Value = "5.04"
$ws.Range("E36").Value = "  +10.78%  "

$ws.Range("E37").Value = "  +3.34%  "

$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").Value = "0.126"
$ws.Range("E39").Value = "  -3.32%  "

$ws.Range("E40").Value = "  -2.89%  "

$ws.Range("E41").Value = "  -5.47%  "

$ws.Range("D42").Value = "71.93"
$ws.Range("E42").Value = "  -2.94%  "

$ws.Range("D43").Value = "0.232"
$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").Value = "12.55"
$ws.Range("E45").Value = "  -9.70%  "

$ws.Range("D46").Value = "1.31"
$ws.Range("E46").Value = "  -4.55%  "

$ws.Range("D47").Value = "5.41"
$ws.Range("E47").Value = "  -6.83%  "

$ws.Range("E48").Value = "  +3.70%  "

$ws.Range("D49").Value = "102.93"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  +7.80%  "

$ws.Range("D51").Value = "8.41"
$ws.Range("E51").Value = "  -1.53%  "
